# Updates the "南宁-漫展信息" workbook (gh-pages generated output) to match
# the scraped data as of commit 456a3b4:
#   - bump a handful of "want to go" (F column) counters on the 展览 and
#     全部类型 sheets
#   - add two newly-scraped rows (Russian ballet shows) to the 演出 sheet
#     and to the combined 全部类型 sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 展览 (Exhibitions) — counter bumps only, no structural change
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 3291   # was 3282
$wsExpo.Range("F6").Value  = 2112   # was 2111
$wsExpo.Range("F10").Value = 1204   # was 1203
$wsExpo.Range("F11").Value = 218    # was 217
$wsExpo.Range("F12").Value = 1226   # was 1214
$wsExpo.Range("F13").Value = 100    # was 99

# ---------------------------------------------------------------------
# 演出 (Performances) — insert the two new ballet rows after row 2
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Rows.Item(3).Insert()
$wsShow.Rows.Item(4).Insert()

$wsShow.Range("A3").Value = 2
$wsShow.Range("B3").NumberFormat = "@"
$wsShow.Range("B3").Value = "2024-07-18"
$wsShow.Range("C3").Value = "南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》"
$wsShow.Range("D3").Value = "龙堤路25号 广西文化艺术中心"
$wsShow.Range("E3").Value = "2024.07.18 20:00-07.18 21:30"
$wsShow.Range("F3").Value = 0
$wsShow.Range("G3").Value = 108
$wsShow.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=85816"
$wsShow.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg"
$wsShow.Range("A3").Font.Bold = $true
$wsShow.Range("A3").HorizontalAlignment = -4108
$wsShow.Range("A3").VerticalAlignment = -4160
$wsShow.Range("A3").Borders.LineStyle = 1

$wsShow.Range("A4").Value = 3
$wsShow.Range("B4").NumberFormat = "@"
$wsShow.Range("B4").Value = "2024-07-19"
$wsShow.Range("C4").Value = "南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 "
$wsShow.Range("D4").Value = "龙堤路25号 广西文化艺术中心"
$wsShow.Range("E4").Value = "2024.07.19 20:00-07.19 22:00"
$wsShow.Range("F4").Value = 0
$wsShow.Range("G4").Value = 108
$wsShow.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=85831"
$wsShow.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg"
$wsShow.Range("A4").Font.Bold = $true
$wsShow.Range("A4").HorizontalAlignment = -4108
$wsShow.Range("A4").VerticalAlignment = -4160
$wsShow.Range("A4").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 全部类型 (All types) — same counter bumps as 展览, plus the same two
# new ballet rows inserted after row 10 (pushing the former rows
# 11-14 down to 13-16)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 3291   # was 3282
$wsAll.Range("F6").Value = 2112   # was 2111

$wsAll.Rows.Item(11).Insert()
$wsAll.Rows.Item(12).Insert()

$wsAll.Range("A11").Value = 10
$wsAll.Range("B11").NumberFormat = "@"
$wsAll.Range("B11").Value = "2024-07-18"
$wsAll.Range("C11").Value = "南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》"
$wsAll.Range("D11").Value = "龙堤路25号 广西文化艺术中心"
$wsAll.Range("E11").Value = "2024.07.18 20:00-07.18 21:30"
$wsAll.Range("F11").Value = 0
$wsAll.Range("G11").Value = 108
$wsAll.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=85816"
$wsAll.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg"
$wsAll.Range("A11").Font.Bold = $true
$wsAll.Range("A11").HorizontalAlignment = -4108
$wsAll.Range("A11").VerticalAlignment = -4160
$wsAll.Range("A11").Borders.LineStyle = 1

$wsAll.Range("A12").Value = 11
$wsAll.Range("B12").NumberFormat = "@"
$wsAll.Range("B12").Value = "2024-07-19"
$wsAll.Range("C12").Value = "南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 "
$wsAll.Range("D12").Value = "龙堤路25号 广西文化艺术中心"
$wsAll.Range("E12").Value = "2024.07.19 20:00-07.19 22:00"
$wsAll.Range("F12").Value = 0
$wsAll.Range("G12").Value = 108
$wsAll.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85831"
$wsAll.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg"
$wsAll.Range("A12").Font.Bold = $true
$wsAll.Range("A12").HorizontalAlignment = -4108
$wsAll.Range("A12").VerticalAlignment = -4160
$wsAll.Range("A12").Borders.LineStyle = 1

# The rows that used to be 11-14 are now 13-16; their F-column "want to
# go" counters also changed in this data refresh.
$wsAll.Range("F13").Value = 1204   # AB动漫游戏嘉年华, was 1203
$wsAll.Range("F14").Value = 218    # 海棠动漫游戏嘉年华, was 217
$wsAll.Range("F15").Value = 1226   # 良牙动漫夏季盛典, was 1214
$wsAll.Range("F16").Value = 100    # 蔚蓝档案only, was 99
